$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '66.549.78'
$cell.ClearFormats()

$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '  -4.06%  '
$cell.ClearFormats()

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.316.34'
$cell.ClearFormats()

$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '  -1.22%  '
$cell.ClearFormats()

$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = '  +0.04%  '
$cell.ClearFormats()

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '572.47'
$cell.ClearFormats()

$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = '  -3.39%  '
$cell.ClearFormats()

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '182.55'
$cell.ClearFormats()

$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '  -5.23%  '
$cell.ClearFormats()

$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '  +0.08%  '
$cell.ClearFormats()

$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = '  -1.12%  '
$cell.ClearFormats()

$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '  -3.21%  '
$cell.ClearFormats()

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '6.63'
$cell.ClearFormats()

$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '  -1.49%  '
$cell.ClearFormats()

$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '  -4.72%  '
$cell.ClearFormats()

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '3.894.96'
$cell.ClearFormats()

$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = '  -1.09%  '
$cell.ClearFormats()

$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '  -0.39%  '
$cell.ClearFormats()

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '27.15'
$cell.ClearFormats()

$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '  -3.86%  '
$cell.ClearFormats()

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '66.626.89'
$cell.ClearFormats()

$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '  -4.03%  '
$cell.ClearFormats()

$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '  -2.65%  '
$cell.ClearFormats()

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '3.344.03'
$cell.ClearFormats()

$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = '  +0.53%  '
$cell.ClearFormats()

$cell = $ws.Range("B18")
$cell.NumberFormat = "@"
$cell.Value = 'Chainlink'
$cell.ClearFormats()

$cell = $ws.Range("C18")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$cell.ClearFormats()

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '13.77'
$cell.ClearFormats()

$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '  +0.12%  '
$cell.ClearFormats()

$cell = $ws.Range("B19")
$cell.NumberFormat = "@"
$cell.Value = 'BitcoinCash'
$cell.ClearFormats()

$cell = $ws.Range("C19")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell.ClearFormats()

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '437.24'
$cell.ClearFormats()

$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = '  +2.37%  '
$cell.ClearFormats()

$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = '  -2.41%  '
$cell.ClearFormats()

$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '  -1.17%  '
$cell.ClearFormats()

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '73.92'
$cell.ClearFormats()

$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '  +0.81%  '
$cell.ClearFormats()

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.ClearFormats()

$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = '  -0.26%  '
$cell.ClearFormats()

$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '  -0.21%  '
$cell.ClearFormats()

$cell = $ws.Range("B25")
$cell.NumberFormat = "@"
$cell.Value = 'PEPE'
$cell.ClearFormats()

$cell = $ws.Range("C25")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$cell.ClearFormats()

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.0000119'
$cell.ClearFormats()

$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '  -2.33%  '
$cell.ClearFormats()

$cell = $ws.Range("B26")
$cell.NumberFormat = "@"
$cell.Value = 'Kaspa'
$cell.ClearFormats()

$cell = $ws.Range("C26")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$cell.ClearFormats()

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '0.193'
$cell.ClearFormats()

$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '  +1.14%  '
$cell.ClearFormats()

$cell = $ws.Range("B27")
$cell.NumberFormat = "@"
$cell.Value = 'InternetComputer(DFINITY)'
$cell.ClearFormats()

$cell = $ws.Range("C27")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell.ClearFormats()

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '9.06'
$cell.ClearFormats()

$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = '  -5.85%  '
$cell.ClearFormats()

$cell = $ws.Range("B28")
$cell.NumberFormat = "@"
$cell.Value = 'Binance-PegBSC-USD'
$cell.ClearFormats()

$cell = $ws.Range("C28")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$cell.ClearFormats()

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()

$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = '  -0.01%  '
$cell.ClearFormats()

$cell = $ws.Range("B29")
$cell.NumberFormat = "@"
$cell.Value = 'PancakeSwap'
$cell.ClearFormats()

$cell = $ws.Range("C29")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$cell.ClearFormats()

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '1.96'
$cell.ClearFormats()

$cell = $ws.Range("E29")
$cell.NumberFormat = "@"
$cell.Value = '  -2.26%  '
$cell.ClearFormats()

$cell = $ws.Range("B30")
$cell.NumberFormat = "@"
$cell.Value = 'EthereumClassic'
$cell.ClearFormats()

$cell = $ws.Range("C30")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$cell.ClearFormats()

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '22.91'
$cell.ClearFormats()

$cell = $ws.Range("E30")
$cell.NumberFormat = "@"
$cell.Value = '  -0.54%  '
$cell.ClearFormats()

$cell = $ws.Range("B31")
$cell.NumberFormat = "@"
$cell.Value = 'NEARProtocol'
$cell.ClearFormats()

$cell = $ws.Range("C31")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$cell.ClearFormats()

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '5.36'
$cell.ClearFormats()

$cell = $ws.Range("E31")
$cell.NumberFormat = "@"
$cell.Value = '  -4.42%  '
$cell.ClearFormats()

$cell = $ws.Range("B32")
$cell.NumberFormat = "@"
$cell.Value = 'USDe'
$cell.ClearFormats()

$cell = $ws.Range("C32")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$cell.ClearFormats()

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.ClearFormats()

$cell = $ws.Range("E32")
$cell.NumberFormat = "@"
$cell.Value = '  +0.09%  '
$cell.ClearFormats()

$cell = $ws.Range("B33")
$cell.NumberFormat = "@"
$cell.Value = 'Aptos'
$cell.ClearFormats()

$cell = $ws.Range("C33")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell.ClearFormats()

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '6.79'
$cell.ClearFormats()

$cell = $ws.Range("E33")
$cell.NumberFormat = "@"
$cell.Value = '  -3.08%  '
$cell.ClearFormats()

$cell = $ws.Range("B34")
$cell.NumberFormat = "@"
$cell.Value = 'Fetch.AI'
$cell.ClearFormats()

$cell = $ws.Range("C34")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell.ClearFormats()

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.23'
$cell.ClearFormats()

$cell = $ws.Range("E34")
$cell.NumberFormat = "@"
$cell.Value = '  -3.98%  '
$cell.ClearFormats()

$cell = $ws.Range("B35")
$cell.NumberFormat = "@"
$cell.Value = 'ImmutableX'
$cell.ClearFormats()

$cell = $ws.Range("C35")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell.ClearFormats()

$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.51'
$cell.ClearFormats()

$cell = $ws.Range("E35")
$cell.NumberFormat = "@"
$cell.Value = '  -0.93%  '
$cell.ClearFormats()

$cell = $ws.Range("B36")
$cell.NumberFormat = "@"
$cell.Value = 'Monero'
$cell.ClearFormats()

$cell = $ws.Range("C36")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$cell.ClearFormats()

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '160.51'
$cell.ClearFormats()

$cell = $ws.Range("E36")
$cell.NumberFormat = "@"
$cell.Value = '  -2.67%  '
$cell.ClearFormats()

$cell = $ws.Range("B37")
$cell.NumberFormat = "@"
$cell.Value = 'EnergySwap'
$cell.ClearFormats()

$cell = $ws.Range("C37")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$cell.ClearFormats()

$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '27.47'
$cell.ClearFormats()

$cell = $ws.Range("E37")
$cell.NumberFormat = "@"
$cell.Value = '  +1.31%  '
$cell.ClearFormats()

$cell = $ws.Range("B38")
$cell.NumberFormat = "@"
$cell.Value = 'Stacks'
$cell.ClearFormats()

$cell = $ws.Range("C38")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell.ClearFormats()

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '1.85'
$cell.ClearFormats()

$cell = $ws.Range("E38")
$cell.NumberFormat = "@"
$cell.Value = '  -3.92%  '
$cell.ClearFormats()

$cell = $ws.Range("B39")
$cell.NumberFormat = "@"
$cell.Value = 'Maker'
$cell.ClearFormats()

$cell = $ws.Range("C39")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$cell.ClearFormats()

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '2.810.18'
$cell.ClearFormats()

$cell = $ws.Range("E39")
$cell.NumberFormat = "@"
$cell.Value = '  +1.91%  '
$cell.ClearFormats()

$cell = $ws.Range("B40")
$cell.NumberFormat = "@"
$cell.Value = 'Mantle'
$cell.ClearFormats()

$cell = $ws.Range("C40")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$cell.ClearFormats()

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.792'
$cell.ClearFormats()

$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '  -2.46%  '
$cell.ClearFormats()

$cell = $ws.Range("B41")
$cell.NumberFormat = "@"
$cell.Value = 'Filecoin'
$cell.ClearFormats()

$cell = $ws.Range("C41")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell.ClearFormats()

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '4.47'
$cell.ClearFormats()

$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '  -2.46%  '
$cell.ClearFormats()

$cell = $ws.Range("B42")
$cell.NumberFormat = "@"
$cell.Value = 'RenderToken'
$cell.ClearFormats()

$cell = $ws.Range("C42")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$cell.ClearFormats()

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '6.25'
$cell.ClearFormats()

$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '  -3.42%  '
$cell.ClearFormats()

$cell = $ws.Range("B43")
$cell.NumberFormat = "@"
$cell.Value = 'Hedera'
$cell.ClearFormats()

$cell = $ws.Range("C43")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$cell.ClearFormats()

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.0678'
$cell.ClearFormats()

$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = '  -1.40%  '
$cell.ClearFormats()

$cell = $ws.Range("B44")
$cell.NumberFormat = "@"
$cell.Value = 'OKB'
$cell.ClearFormats()

$cell = $ws.Range("C44")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$cell.ClearFormats()

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '40.16'
$cell.ClearFormats()

$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '  -2.49%  '
$cell.ClearFormats()

$cell = $ws.Range("B45")
$cell.NumberFormat = "@"
$cell.Value = 'InjectiveProtocol'
$cell.ClearFormats()

$cell = $ws.Range("C45")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$cell.ClearFormats()

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '24.22'
$cell.ClearFormats()

$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '  -5.00%  '
$cell.ClearFormats()

$cell = $ws.Range("B46")
$cell.NumberFormat = "@"
$cell.Value = 'dogwifhat'
$cell.ClearFormats()

$cell = $ws.Range("C46")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell.ClearFormats()

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.35'
$cell.ClearFormats()

$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '  -6.77%  '
$cell.ClearFormats()

$cell = $ws.Range("B47")
$cell.NumberFormat = "@"
$cell.Value = 'Bittensor'
$cell.ClearFormats()

$cell = $ws.Range("C47")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$cell.ClearFormats()

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '319.64'
$cell.ClearFormats()

$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '  -7.14%  '
$cell.ClearFormats()

$cell = $ws.Range("B48")
$cell.NumberFormat = "@"
$cell.Value = 'VeChain'
$cell.ClearFormats()

$cell = $ws.Range("C48")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell.ClearFormats()

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.0274'
$cell.ClearFormats()

$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = '  -2.97%  '
$cell.ClearFormats()

$cell = $ws.Range("B49")
$cell.NumberFormat = "@"
$cell.Value = 'ONDO'
$cell.ClearFormats()

$cell = $ws.Range("C49")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$cell.ClearFormats()

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.983'
$cell.ClearFormats()

$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '  -2.33%  '
$cell.ClearFormats()

$cell = $ws.Range("B50")
$cell.NumberFormat = "@"
$cell.Value = 'Cosmos'
$cell.ClearFormats()

$cell = $ws.Range("C50")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell.ClearFormats()

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '6.18'
$cell.ClearFormats()

$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = '  -1.64%  '
$cell.ClearFormats()

$cell = $ws.Range("B51")
$cell.NumberFormat = "@"
$cell.Value = 'Stellar'
$cell.ClearFormats()

$cell = $ws.Range("C51")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell.ClearFormats()

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '0.0997'
$cell.ClearFormats()

$cell = $ws.Range("E51")
$cell.NumberFormat = "@"
$cell.Value = '  -1.65%  '
$cell.ClearFormats()
